$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 105
$ws.Range("I2").Value = 275
$ws.Range("J2").Value = 1104
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 296
$ws.Range("M2").Value = 18
$ws.Range("N2").Value = 201
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 6
$ws.Range("R2").Value = 15
$ws.Range("S2").Value = 116
$ws.Range("T2").Value = 204
$ws.Range("V2").Value = 1739
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1641
$ws.Range("Z2").Value = 29
$ws.Range("AA2").Value = 8
